$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 236.66667
$ws.Range("I31").Value = 236.66667
$ws.Range("K31").Value = 710.00001
$ws.Range("M31").Value = -480.00001
$ws.Range("H38").Value = 258.55554
$ws.Range("I38").Value = 165.875
$ws.Range("K38").Value = 497.625
$ws.Range("M38").Value = -125.625
$ws.Range("H39").Value = 768.4737
$ws.Range("I39").Value = 912.55554
$ws.Range("J39").Value = 638.8
$ws.Range("K39").Value = 2737.66662
$ws.Range("L39").Value = 1916.4
$ws.Range("M39").Value = -2441.66662
$ws.Range("N39").Value = -2508.4
$ws.Range("H107").Value = 1659.625
$ws.Range("J107").Value = 3181
$ws.Range("L107").Value = 3181
$ws.Range("N107").Value = -7021
$ws.Range("H138").Value = 2173.01
$ws.Range("I138").Value = 1666.8889
$ws.Range("J138").Value = 2223.066
$ws.Range("K138").Value = 5000.6667
$ws.Range("L138").Value = 6669.197999999999
$ws.Range("M138").Value = 139.3333000000002
$ws.Range("N138").Value = -16949.198
$ws.Range("H139").Value = 44488
$ws.Range("J139").Value = 46653.332
$ws.Range("L139").Value = 46653.332
$ws.Range("N139").Value = -56933.332
$ws.Range("H140").Value = 45383.168
$ws.Range("J140").Value = 45383.168
$ws.Range("L140").Value = 45383.168
$ws.Range("N140").Value = -55743.168
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1011.7
$ws.Range("I61").Value = 848.4286
$ws.Range("J61").Value = 1392.6666
$ws.Range("K61").Value = 848.4286
$ws.Range("L61").Value = 1392.6666
$ws.Range("M61").Value = -636.4286
$ws.Range("N61").Value = -1816.6666
$ws.Range("H110").Value = 1298.2
$ws.Range("I110").Value = 1095.2
$ws.Range("K110").Value = 1095.2
$ws.Range("M110").Value = 949.8
$ws.Range("H136").Value = 1011.7
$ws.Range("I136").Value = 848.4286
$ws.Range("J136").Value = 1392.6666
$ws.Range("K136").Value = 2545.2858
$ws.Range("L136").Value = 4177.9998
$ws.Range("M136").Value = 4.714200000000346
$ws.Range("N136").Value = -9277.9998
$ws.Range("H139").Value = 50214.5
$ws.Range("J139").Value = 50214.5
$ws.Range("L139").Value = 50214.5
$ws.Range("N139").Value = -60494.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 24695.223
$ws.Range("I82").Value = 18651.4
$ws.Range("K82").Value = 18651.4
$ws.Range("M82").Value = -18268.4
$ws.Range("H85").Value = 24695.223
$ws.Range("I85").Value = 18651.4
$ws.Range("K85").Value = 18651.4
$ws.Range("M85").Value = -17325.4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1375.6842
$ws.Range("I99").Value = 1469.4445
$ws.Range("K99").Value = 1469.4445
$ws.Range("M99").Value = 28.55549999999994
$ws.Range("H105").Value = 703
$ws.Range("I105").Value = 553.3333
$ws.Range("J105").Value = 759.125
$ws.Range("K105").Value = 553.3333
$ws.Range("L105").Value = 759.125
$ws.Range("M105").Value = 1193.6667
$ws.Range("N105").Value = -4253.125
$ws.Range("H126").Value = 1375.6842
$ws.Range("I126").Value = 1469.4445
$ws.Range("K126").Value = 4408.333500000001
$ws.Range("M126").Value = -1938.333500000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1653.6538
$ws.Range("I5").Value = 1404.3334
$ws.Range("J5").Value = 2700.8
$ws.Range("K5").Value = 4213.0002
$ws.Range("L5").Value = 8102.400000000001
$ws.Range("M5").Value = -4101.0002
$ws.Range("N5").Value = -8326.400000000001
$ws.Range("H63").Value = 5497.75
$ws.Range("I63").Value = 996
$ws.Range("J63").Value = 6998.3335
$ws.Range("K63").Value = 2988
$ws.Range("L63").Value = 20995.0005
$ws.Range("M63").Value = -2239
$ws.Range("N63").Value = -22493.0005
$ws.Range("H66").Value = 5497.75
$ws.Range("I66").Value = 996
$ws.Range("J66").Value = 6998.3335
$ws.Range("K66").Value = 8964
$ws.Range("L66").Value = 62985.0015
$ws.Range("M66").Value = -5220
$ws.Range("N66").Value = -70473.0015
$ws.Range("H68").Value = 1661.5
$ws.Range("J68").Value = 1860.2307
$ws.Range("L68").Value = 5580.6921
$ws.Range("N68").Value = -7202.6921
$ws.Range("H70").Value = 4852.9414
$ws.Range("J70").Value = 5471.4287
$ws.Range("L70").Value = 16414.2861
$ws.Range("N70").Value = -17044.2861
$ws.Range("H71").Value = 1661.5
$ws.Range("J71").Value = 1860.2307
$ws.Range("L71").Value = 16742.0763
$ws.Range("N71").Value = -24854.0763
$ws.Range("H73").Value = 4852.9414
$ws.Range("J73").Value = 5471.4287
$ws.Range("L73").Value = 16414.2861
$ws.Range("N73").Value = -18598.2861
$ws.Range("H121").Value = 638.8
$ws.Range("I121").Value = 297
$ws.Range("J121").Value = 866.6667
$ws.Range("K121").Value = 891
$ws.Range("L121").Value = 2600.0001
$ws.Range("M121").Value = 419
$ws.Range("N121").Value = -5220.0001
$ws.Range("H135").Value = 1653.6538
$ws.Range("I135").Value = 1404.3334
$ws.Range("J135").Value = 2700.8
$ws.Range("K135").Value = 12639.0006
$ws.Range("L135").Value = 24307.2
$ws.Range("M135").Value = -10104.0006
$ws.Range("N135").Value = -29377.2
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 30000
$ws.Range("J93").Value = 30000
$ws.Range("L93").Value = 30000
$ws.Range("N93").Value = -33744
$ws.Range("H132").Value = 2592.4443
$ws.Range("I132").Value = 2213.889
$ws.Range("J132").Value = 3349.5557
$ws.Range("K132").Value = 6641.667
$ws.Range("L132").Value = 10048.6671
$ws.Range("M132").Value = -4111.667
$ws.Range("N132").Value = -15108.6671
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H98").Value = 35355
$ws.Range("J98").Value = 35355
$ws.Range("L98").Value = 35355
$ws.Range("N98").Value = -41345
$ws.Range("H122").Value = 17860926
$ws.Range("I122").Value = 62503130
$ws.Range("J122").Value = 4046.2
$ws.Range("K122").Value = 187509390
$ws.Range("L122").Value = 12138.6
$ws.Range("M122").Value = -187506940
$ws.Range("N122").Value = -17038.6
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H95").Value = 29499.5
$ws.Range("J95").Value = 29499.5
$ws.Range("L95").Value = 29499.5
$ws.Range("N95").Value = -34991.5
$ws.Range("H132").Value = 4538.0835
$ws.Range("I132").Value = 5766.5713
$ws.Range("J132").Value = 2818.2
$ws.Range("K132").Value = 17299.7139
$ws.Range("L132").Value = 8454.599999999999
$ws.Range("M132").Value = -14769.7139
$ws.Range("N132").Value = -13514.6
